# Weekly data refresh: insert the latest week's two new price records
# ("Inferno" Primera and Segunda, dated 2022-10-03 / serial 44837) at the
# top of the data block (row 80), pushing every existing record down by
# two rows (old row 80 -> new row 82, ..., old row 104 -> new row 106).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room: insert two blank rows at row 80. Doing this twice (rather than
# a 2-row range insert) mirrors inserting them one at a time and reliably
# shifts all data below down by two rows while inheriting row 79/80's
# formatting (needed for the date-formatted column D).
$ws.Rows.Item(80).Insert()
$ws.Rows.Item(80).Insert()

# New row 80: Ají / Inferno / Primera
$ws.Range("A80").Value = 1
$ws.Range("B80").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C80").Value = "Arica y Parinacota"
$ws.Range("D80").Value = 44837
$ws.Range("E80").Value = 15
$ws.Range("F80").Value = 100112021
$ws.Range("G80").Value = "Ají"
$ws.Range("H80").Value = "Inferno"
$ws.Range("I80").Value = "Primera"
$ws.Range("J80").Value = 120
$ws.Range("K80").Value = 17000
$ws.Range("L80").Value = 18000
$ws.Range("M80").Value = 17500
$ws.Range("N80").Value = "$/caja 15 kilos"
$ws.Range("O80").Value = "Región de Arica y Parinacota"
$ws.Range("P80").Value = 1167
$ws.Range("Q80").Value = 15
$ws.Range("R80").Value = "Hortaliza"

# New row 81: Ají / Inferno / Segunda
$ws.Range("A81").Value = 1
$ws.Range("B81").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C81").Value = "Arica y Parinacota"
$ws.Range("D81").Value = 44837
$ws.Range("E81").Value = 15
$ws.Range("F81").Value = 100112021
$ws.Range("G81").Value = "Ají"
$ws.Range("H81").Value = "Inferno"
$ws.Range("I81").Value = "Segunda"
$ws.Range("J81").Value = 160
$ws.Range("K81").Value = 13000
$ws.Range("L81").Value = 14000
$ws.Range("M81").Value = 13500
$ws.Range("N81").Value = "$/caja 15 kilos"
$ws.Range("O81").Value = "Región de Arica y Parinacota"
$ws.Range("P81").Value = 900
$ws.Range("Q81").Value = 15
$ws.Range("R81").Value = "Hortaliza"
